$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same
# three rows (2, 8, 9) whose "想去人数" (F column) counts were updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 599
    $ws.Range("F8").Value = 567
    $ws.Range("F9").Value = 3750
}
